$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.941.86'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '1.620.63'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.63'
$ws.Range('E5').Value = '  -0.99%  '
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.491'
$ws.Range('E7').Value = '  -2.86%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0622'
$ws.Range('E8').Value = '  -0.89%  '
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('E10').Value = '  -2.08%  '
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = '1.845.83'
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('D13').Value = '1.617.90'
$ws.Range('E13').Value = '  -4.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.15'
$ws.Range('E14').Value = '  -1.64%  '
$ws.Range('E15').Value = '  -1.84%  '
$ws.Range('D16').Value = '25.961.81'
$ws.Range('E16').Value = '  -0.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.73'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('E18').Value = '  -1.49%  '
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '191.93'
$ws.Range('E20').Value = '  +0.44%  '
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('E22').Value = '  -0.99%  '
$ws.Range('E23').Value = '  -2.11%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.40'
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('E26').Value = '  -0.49%  '
$ws.Range('E27').Value = '  -4.03%  '
$ws.Range('E28').Value = '  -1.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.22'
$ws.Range('E29').Value = '  -0.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.22'
$ws.Range('E30').Value = '  -1.19%  '
$ws.Range('E31').Value = '  -1.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.12'
$ws.Range('E32').Value = '  -1.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.41'
$ws.Range('E34').Value = '  -0.90%  '
$ws.Range('E35').Value = '  -1.17%  '
$ws.Range('D36').Value = '1.128.46'
$ws.Range('E36').Value = '  -0.34%  '
$ws.Range('E37').Value = '  -3.98%  '
$ws.Range('E38').Value = '  -1.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.519'
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '97.91'
$ws.Range('E41').Value = '  -1.09%  '
$ws.Range('D42').Value = '1.757.18'
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.756'
$ws.Range('E43').Value = '  -4.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.13'
$ws.Range('E44').Value = '  -3.41%  '
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.08'
$ws.Range('E47').Value = '  -2.68%  '
$ws.Range('E48').Value = '  -1.83%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.01'
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.46'
$ws.Range('E51').Value = '  -1.43%  '
